$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 2.147570471799392
$arr[0,1] = -1.292459514723937
$arr[0,2] = -1.353633693138139
$arr[0,3] = 0.6878222440181159
$arr[0,4] = 0.1052277069844749
$arr[0,5] = 0.2009027477938327
$arr[0,6] = -0.0595162892048901
$arr[0,7] = 0.7759387069536011
$arr[0,8] = 0.6758593470509796
$arr[0,9] = 0.3220728409712834
$ws.Range("B2:K2").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -3.44002998652333
$arr[0,1] = -3.501204164937531
$arr[0,2] = -1.459748227781277
$arr[0,3] = -2.042342764814918
$arr[0,4] = -1.94666772400556
$arr[0,5] = -2.207086761004283
$arr[0,6] = -1.371631764845791
$arr[0,7] = -1.471711124748413
$arr[0,8] = -1.825497630828109
$arr[0,9] = -1.677320240395664
$ws.Range("B3:K3").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.06117417841420103
$arr[0,1] = 1.980281758742053
$arr[0,2] = 1.397687221708412
$arr[0,3] = 1.49336226251777
$arr[0,4] = 1.232943225519047
$arr[0,5] = 2.068398221677539
$arr[0,6] = 1.968318861774917
$arr[0,7] = 1.614532355695221
$arr[0,8] = 1.762709746127666
$arr[0,9] = 1.920155735131821
$ws.Range("B4:K4").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 2.041455937156254
$arr[0,1] = 1.458861400122613
$arr[0,2] = 1.554536440931971
$arr[0,3] = 1.294117403933248
$arr[0,4] = 2.12957240009174
$arr[0,5] = 2.029493040189118
$arr[0,6] = 1.675706534109422
$arr[0,7] = 1.823883924541867
$arr[0,8] = 1.981329913546022
$arr[0,9] = 1.480670457019951
$ws.Range("B5:K5").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.5825945370336409
$arr[0,1] = -0.4869194962242832
$arr[0,2] = -0.747338533223006
$arr[0,3] = 0.0881164629354852
$arr[0,4] = -0.01196289696713632
$arr[0,5] = -0.3657494030468326
$arr[0,6] = -0.2175720126143872
$arr[0,7] = -0.06012602361023223
$arr[0,8] = -0.560785480136303
$arr[0,9] = -0.2804275996008339
$ws.Range("B6:K6").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.09567504080935779
$arr[0,1] = -0.164743996189365
$arr[0,2] = 0.6707109999691262
$arr[0,3] = 0.5706316400665047
$arr[0,4] = 0.2168451339868084
$arr[0,5] = 0.3650225244192538
$arr[0,6] = 0.5224685134234088
$arr[0,7] = 0.02180905689733798
$arr[0,8] = 0.3021669374328071
$arr[0,9] = 0.1301920795877955
$ws.Range("B7:K7").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.2604190369987228
$arr[0,1] = 0.5750359591597685
$arr[0,2] = 0.4749565992571469
$arr[0,3] = 0.1211700931774507
$arr[0,4] = 0.269347483609896
$arr[0,5] = 0.426793472614051
$arr[0,6] = -0.07386598391201982
$arr[0,7] = 0.2064918966234494
$arr[0,8] = 0.0345170387784377
$arr[0,9] = 0.3403795785247692
$ws.Range("B8:K8").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.8354549961584912
$arr[0,1] = 0.7353756362558697
$arr[0,2] = 0.3815891301761735
$arr[0,3] = 0.5297665206086188
$arr[0,4] = 0.6872125096127738
$arr[0,5] = 0.186553053086703
$arr[0,6] = 0.4669109336221722
$arr[0,7] = 0.2949360757771605
$arr[0,8] = 0.600798615523492
$arr[0,9] = -0.01153133716379418
$ws.Range("B9:K9").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.1000793599026215
$arr[0,1] = -0.4538658659823178
$arr[0,2] = -0.3056884755498724
$arr[0,3] = -0.1482424865457174
$arr[0,4] = -0.6489019430717882
$arr[0,5] = -0.3685440625363191
$arr[0,6] = -0.5405189203813308
$arr[0,7] = -0.2346563806349992
$arr[0,8] = -0.8469863333222853
$arr[0,9] = -0.1586149481231739
$ws.Range("B10:K10").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.3537865060796963
$arr[0,1] = -0.2056091156472509
$arr[0,2] = -0.04816312664309591
$arr[0,3] = -0.5488225831691667
$arr[0,4] = -0.2684647026336975
$arr[0,5] = -0.4404395604787092
$arr[0,6] = -0.1345770207323777
$arr[0,7] = -0.7469069734196638
$arr[0,8] = -0.05853558822055238
$arr[0,9] = -0.3352267436446591
$ws.Range("B11:K11").Value2 = $arr

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 0.1481773904324453
$arr[0,1] = 0.3056233794366003
$arr[0,2] = -0.1950360770894705
$arr[0,3] = 0.08532180344599868
$arr[0,4] = -0.08665305439901295
$arr[0,5] = 0.2192094853473185
$arr[0,6] = -0.3931204673399676
$arr[0,7] = 0.2952509178591439
$arr[0,8] = 0.01855976243503714
$ws.Range("B12:J12").Value2 = $arr

$arr = New-Object "object[,]" 1,8
$arr[0,0] = 0.157445989004155
$arr[0,1] = -0.3432134675219158
$arr[0,2] = -0.06285558698644665
$arr[0,3] = -0.2348304448314583
$arr[0,4] = 0.0710320949148732
$arr[0,5] = -0.541297857772413
$arr[0,6] = 0.1470735274266985
$arr[0,7] = -0.1296176279974082
$ws.Range("B13:I13").Value2 = $arr

$arr = New-Object "object[,]" 1,7
$arr[0,0] = -0.5006594565260708
$arr[0,1] = -0.2203015759906016
$arr[0,2] = -0.3922764338356133
$arr[0,3] = -0.0864138940892818
$arr[0,4] = -0.698743846776568
$arr[0,5] = -0.01037246157745647
$arr[0,6] = -0.2870636170015632
$ws.Range("B14:H14").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = 0.2803578805354692
$arr[0,1] = 0.1083830226904575
$arr[0,2] = 0.414245562436789
$arr[0,3] = -0.1980843902504972
$arr[0,4] = 0.4902869949486143
$arr[0,5] = 0.2135958395245076
$ws.Range("B15:G15").Value2 = $arr

$arr = New-Object "object[,]" 1,5
$arr[0,0] = -0.1719748578450117
$arr[0,1] = 0.1338876819013198
$arr[0,2] = -0.4784422707859664
$arr[0,3] = 0.2099291144131452
$arr[0,4] = -0.06676204101096155
$ws.Range("B16:F16").Value2 = $arr

$arr = New-Object "object[,]" 1,4
$arr[0,0] = 0.3058625397463315
$arr[0,1] = -0.3064674129409547
$arr[0,2] = 0.3819039722581568
$arr[0,3] = 0.1052128168340501
$ws.Range("B17:E17").Value2 = $arr

$arr = New-Object "object[,]" 1,3
$arr[0,0] = -0.6123299526872862
$arr[0,1] = 0.07604143251182532
$arr[0,2] = -0.2006497229122814
$ws.Range("B18:D18").Value2 = $arr

$arr = New-Object "object[,]" 1,2
$arr[0,0] = 0.6883713851991116
$arr[0,1] = 0.4116802297750048
$ws.Range("B19:C19").Value2 = $arr

$arr = New-Object "object[,]" 1,1
$arr[0,0] = -0.2766911554241067
$ws.Range("B20:B20").Value2 = $arr
